$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the paragraphs we need to touch by their visible text content.
# ---------------------------------------------------------------------

$week2Idx = -1
$todoWeek2Idx = -1
$reactNativeIdx = -1
$androidEmuIdx = -1
$structuurIdx = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    $trimmed = $t.TrimEnd()

    if ($week2Idx -eq -1 -and $t -like "Week 2*") {
        $week2Idx = $i
    }
    if ($week2Idx -ne -1 -and $todoWeek2Idx -eq -1 -and $i -gt $week2Idx -and $trimmed -eq "To-do") {
        $todoWeek2Idx = $i
    }
    if ($t -like "*beginnen programmeren*") {
        $reactNativeIdx = $i
    }
    if ($t -like "*native installeren + instellen*") {
        $androidEmuIdx = $i
    }
    if ($t -like "Structuur masterproef af*") {
        $structuurIdx = $i
    }
}

$paraTodoWeek2       = $d.Paragraphs.Item($todoWeek2Idx)
$paraReactNative     = $d.Paragraphs.Item($reactNativeIdx)
$paraAndroidEmulator = $d.Paragraphs.Item($androidEmuIdx)
$paraStructuur       = $d.Paragraphs.Item($structuurIdx)

# Known-good "space run" + "checkmark run" pair to clone formatting from
# (identical rPr is used throughout the document for these markers).
$checkSrcEnd = $paraStructuur.Range.End - 1
$checkSrc = $d.Range($checkSrcEnd - 1, $checkSrcEnd)
$spaceSrc = $d.Range($checkSrcEnd - 2, $checkSrcEnd - 1)

# ---------------------------------------------------------------------
# 1) Append a checkmark run right after "Native app beginnen programmeren"
# ---------------------------------------------------------------------
$insertPos = $paraReactNative.Range.End - 1
$dest = $d.Range($insertPos, $insertPos)
$dest.FormattedText = $checkSrc.FormattedText

# ---------------------------------------------------------------------
# 2) Append a space run + checkmark run after "...instellen", replacing
#    the bookmark that used to sit there.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$insertPos2 = $paraAndroidEmulator.Range.End - 1
$destSpace = $d.Range($insertPos2, $insertPos2)
$destSpace.FormattedText = $spaceSrc.FormattedText

$insertPos3 = $paraAndroidEmulator.Range.End - 1
$destCheck = $d.Range($insertPos3, $insertPos3)
$destCheck.FormattedText = $checkSrc.FormattedText

# ---------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark at the end of the Week 2 "To-do"
#    paragraph. Bookmarking a collapsed range that sits exactly at a
#    paragraph boundary needs a small workaround: insert a throw-away
#    character after the target point, bookmark in front of it, then
#    remove the throw-away character again.
# ---------------------------------------------------------------------
$goBackPos = $paraTodoWeek2.Range.End - 1
$tmp = $d.Range($goBackPos, $goBackPos)
$tmp.InsertAfter("X")
$bmRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($goBackPos, $goBackPos + 1).Delete()

Write-Output "done"
